# Updates on COvid-19 in Kenya 24 April 2020
# Adds a new data row (row 42) to the Kenya COVID-19 tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 42

# A: Date  (24-Apr-2020 => serial 43945), formatted/styled like the preceding date cells
$ws.Cells.Item($newRow, 1).Value = 43945
$ws.Cells.Item($newRow, 1).NumberFormat = "d-mmm-yy"
$ws.Cells.Item($newRow, 1).HorizontalAlignment = -4108

# B: New Cases
$ws.Cells.Item($newRow, 2).Value = 16

# C: Tested
$ws.Cells.Item($newRow, 3).Value = 946

# D: Travelled From
$ws.Cells.Item($newRow, 4).Value = "None"

# E: County
$ws.Cells.Item($newRow, 5).Value = "Mombasa(5),Nairobi(11)"

# F: Aggregation
$ws.Cells.Item($newRow, 6).Value = 336

# G: Case Type
$ws.Cells.Item($newRow, 7).Value = "Community(16)"

# H: Recover
$ws.Cells.Item($newRow, 8).Value = 5

# I: Death
$ws.Cells.Item($newRow, 9).Value = 0

# Update the visible view/selection to reflect the new active cell location
$ws.Application.ActiveWindow.ScrollRow = 20
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("I42").Select()
